$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: copy date formatting from A27 into A28, then fill values ---
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A28").Value = 43817
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 1
$ws.Range("M28").Value = 3

# --- Row 29: copy date formatting from A28 into A29, then fill values ---
$ws.Range("A28").Copy()
$ws.Range("A29").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A29").Value = 43818
$ws.Range("J29").Value = 41

# --- Update the selection to match the new active cell ---
$ws.Range("B29").Select() | Out-Null
